$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bibi Cell Mundi
$ws.Range("L2").Value = 16641.22
$ws.Range("M2").Value = 1882
$ws.Range("AG2").Value = 106227.95

# Row 3 - Bibi Cell Vieiralves
$ws.Range("L3").Value = 2865
$ws.Range("M3").Value = 3923
$ws.Range("AG3").Value = 44714.05

# Row 4 - Bibi Cell Manauara
$ws.Range("L4").Value = 5175
$ws.Range("M4").Value = 2532
$ws.Range("N4").Value = 1811
$ws.Range("AG4").Value = 41272.4

# Row 5 - Bibi Cell Ponta Negra
$ws.Range("L5").Value = 717
$ws.Range("M5").Value = 1708.9
$ws.Range("N5").Value = 1178.9
$ws.Range("AG5").Value = 33221.85

# Row 6 - total
$ws.Range("L6").Value = 25398.22
$ws.Range("M6").Value = 10045.9
$ws.Range("N6").Value = 2989.9
$ws.Range("AG6").Value = 225436.25
